# Applies the "Adjust costs where relying on US data" change:
#  - About sheet: add a note + EU/US reference values (rows 26-28)
#  - ICtPSFfL sheet: multiply row 7 (incremental cost) by the EU:US ratio

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A26").Value = "We adjust for the EU data by the ratio of EU:US pre-tax transportation biofuel costs (see file fuels/BFPaT for the EU and US models)."

$wsAbout.Range("A27").Value = "EU"
$wsAbout.Range("B27").Value = 0.000018152570386688024

$wsAbout.Range("A28").Value = "US"
$wsAbout.Range("B28").Value = 0.000012337034592036476

$wsFuel = $wb.Worksheets.Item("ICtPSFfL")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK")

foreach ($col in $cols) {
    $addr = "$col" + "7"
    $wsFuel.Range($addr).Formula = "=MAX(Calcs!" + $col + "35,0)*(About!`$B`$27/About!`$B`$28)"
}
